$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '22.482.41'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.573.16'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -0.13%  '
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '292.09'
$ws.Range("E6").Value = '  +0.11%  '
$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '0.3723'
$ws.Range("E7").Value = '  -1.14%  '
$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '49.81'
$ws.Range("E8").Value = '  -0.03%  '
$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '0.3402'
$ws.Range("E9").Value = '  -0.48%  '
$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '1.149'
$ws.Range("E10").Value = '  +0.05%  '
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '0.07551'
$ws.Range("E11").Value = '  -1.14%  '
$ws.Range("E12").Value = '  -0.05%  '
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '21.27'
$ws.Range("E13").Value = '  +0.27%  '
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '6.044'
$ws.Range("E14").Value = '  +0.54%  '
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '6.968'
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").Value = '1.573.16'
$ws.Range("E16").Value = '  -0.21%  '
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '0.00001126'
$ws.Range("E17").Value = '  -0.80%  '
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '90.85'
$ws.Range("E18").Value = '  +0.53%  '
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '0.06770'
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("E20").Value = '  -0.07%  '
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '6.307'
$ws.Range("E21").Value = '  +1.45%  '
$ws.Range("E22").Value = '  -2.19%  '
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '12.17'
$ws.Range("E23").Value = '  +1.20%  '
$ws.Range("E24").Value = '  +0.28%  '
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '2.371'
$ws.Range("E25").Value = '  -0.99%  '
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '2.630'
$ws.Range("E26").Value = '  -1.41%  '
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '20.02'
$ws.Range("E27").Value = '  -0.59%  '
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '149.47'
$ws.Range("E28").Value = '  +1.58%  '
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '5.052'
$ws.Range("E29").Value = '  +0.33%  '
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '125.60'
$ws.Range("E30").Value = '  -0.83%  '
$ws.Range("D31").Value = '1.746.45'
$ws.Range("E31").Value = '  -0.22%  '
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '1.081'
$ws.Range("E32").Value = '  +9.80%  '
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '6.210'
$ws.Range("E33").Value = '  +1.14%  '
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '2.016'
$ws.Range("E34").Value = '  +0.31%  '
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '9.815'
$ws.Range("E35").Value = '  -3.64%  '
$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '0.08361'
$ws.Range("E36").Value = '  -1.76%  '
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '0.02478'
$ws.Range("E37").Value = '  -2.47%  '
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '0.2302'
$ws.Range("E38").Value = '  -0.57%  '
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '1.339'
$ws.Range("E39").Value = '  -2.09%  '
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '0.06542'
$ws.Range("E40").Value = '  +0.27%  '
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '5.446'
$ws.Range("E41").Value = '  +0.05%  '
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '11.36'
$ws.Range("E42").Value = '  -0.42%  '
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '0.6243'
$ws.Range("E43").Value = '  -2.07%  '
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '14.14'
$ws.Range("E44").Value = '  +0.71%  '
$ws.Range("E45").Value = '  -0.07%  '
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '3.815'
$ws.Range("E46").Value = '  +0.58%  '
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '0.5849'
$ws.Range("E47").Value = '  -2.05%  '
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '130.47'
$ws.Range("E48").Value = '  +4.59%  '
$ws.Range("E49").Value = '  -1.02%  '
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '1.214'
$ws.Range("E50").Value = '  -5.36%  '
$ws.Range("E51").Value = '  +0.12%  '

# Clear the temporary text-number-format so cells keep default (unstyled) appearance
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D50").ClearFormats()
